$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..3) {
    $ws.Cells.Item($r, 4).Value = -0.0246          # D
    $ws.Cells.Item($r, 5).ClearContents()          # E removed
    $ws.Cells.Item($r, 7).Value = 0.08467966573816156   # G
    $ws.Cells.Item($r, 8).Value = 0.08467966573816156   # H
    $ws.Cells.Item($r, 9).Value = -0.02116991643454039  # I
    $ws.Cells.Item($r, 10).Value = -0.02116991643454039 # J
    $ws.Cells.Item($r, 11).Value = -0.083                # K
    $ws.Cells.Item($r, 12).Value = -0.02311977715877438  # L
    $ws.Cells.Item($r, 13).Value = 0.212                 # M
    $ws.Cells.Item($r, 14).Value = 0.0157037037037037   # N
    $ws.Cells.Item($r, 15).Value = -2.554216867469879   # O
    $ws.Cells.Item($r, 16).Value = 0.212                 # P
    $ws.Cells.Item($r, 17).Value = 0.0157037037037037   # Q
    $ws.Cells.Item($r, 18).Value = -2.554216867469879   # R
    $ws.Cells.Item($r, 21).Value = 0.509                 # U
    $ws.Cells.Item($r, 22).Value = 0.0377037037037037   # V
    $ws.Cells.Item($r, 23).Value = -0.007280701754385965 # W
    $ws.Cells.Item($r, 24).Value = 0.06322667656689483  # X
    $ws.Cells.Item($r, 25).Value = -0.0705073783212808  # Y
    $ws.Cells.Item($r, 26).Value = 0.3374377291098787   # Z
    $ws.Cells.Item($r, 27).Value = -0.007143528527117209 # AA
    $ws.Cells.Item($r, 28).Value = 0.06322667656689483  # AB
    $ws.Cells.Item($r, 29).Value = -0.07037020509401204 # AC
    $ws.Cells.Item($r, 33).Value = -0.509                # AG
    $ws.Cells.Item($r, 36).Value = -0.03918097144176738 # AJ
    $ws.Cells.Item($r, 37).Value = -0.04851777714231246 # AK
    $ws.Cells.Item($r, 39).Value = -0.064                # AM
    $ws.Cells.Item($r, 42).Value = -3.161490683229814   # AP
    $ws.Cells.Item($r, 43).Value = 1.1875                # AQ
}
